$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F ("p_value"), matching the bold/centered header
# style already used for A1:E1.
$ws.Cells.Item(1, 6).Value = "p_value"
$ws.Cells.Item(1, 6).Font.Bold = $true
$ws.Cells.Item(1, 6).HorizontalAlignment = -4108

# p-values for the permutation test, one per data row (rows 2-251).
$pValues = @(
    0.8118811881188119,
    0.06930693069306931,
    0.495049504950495,
    0.9207920792079208,
    0.3663366336633663,
    0.8613861386138614,
    0.3168316831683168,
    0.2574257425742574,
    0.2871287128712871,
    0.1683168316831683,
    0.5841584158415841,
    0.2475247524752475,
    0.7326732673267327,
    0.198019801980198,
    0.693069306930693,
    0.3861386138613861,
    0.3267326732673267,
    0.7227722772277227,
    0.8316831683168316,
    0.4752475247524752,
    0.7920792079207921,
    0.4554455445544555,
    0.1287128712871287,
    0.2772277227722773,
    0.2772277227722773,
    0.09900990099009901,
    0.405940594059406,
    0.5742574257425742,
    0.8217821782178217,
    0.0297029702970297,
    0.3069306930693069,
    0.7722772277227723,
    0.1683168316831683,
    0.1485148514851485,
    0.7623762376237624,
    0.5544554455445545,
    0.5247524752475248,
    0.8712871287128713,
    0.8910891089108911,
    0.693069306930693,
    0.8118811881188119,
    0.3069306930693069,
    0.1089108910891089,
    0.0891089108910891,
    0.7227722772277227,
    0.3366336633663367,
    0.7920792079207921,
    0.0198019801980198,
    0.5247524752475248,
    0.2277227722772277,
    1.0,
    0.5148514851485149,
    0.1485148514851485,
    0.2871287128712871,
    0.07920792079207921,
    0.5643564356435643,
    0.07920792079207921,
    0.5346534653465347,
    0.7524752475247525,
    0.4455445544554456,
    0.900990099009901,
    0.1881188118811881,
    0.6435643564356436,
    0.1188118811881188,
    0.7722772277227723,
    0.0297029702970297,
    0.7326732673267327,
    0.2673267326732673,
    0.04950495049504951,
    0.3960396039603961,
    0.3069306930693069,
    0.5643564356435643,
    0.7524752475247525,
    0.1485148514851485,
    0.504950495049505,
    0.2475247524752475,
    0.2079207920792079,
    0.3366336633663367,
    0.3663366336633663,
    0.0891089108910891,
    0.5841584158415841,
    0.1782178217821782,
    0.0891089108910891,
    0.2079207920792079,
    0.4554455445544555,
    0.5148514851485149,
    0.6534653465346535,
    0.7326732673267327,
    0.1089108910891089,
    0.7029702970297029,
    0.7029702970297029,
    0.6039603960396039,
    0.9603960396039604,
    0.0297029702970297,
    0.1287128712871287,
    0.3663366336633663,
    0.4752475247524752,
    0.8514851485148515,
    0.4554455445544555,
    0.04950495049504951,
    0.297029702970297,
    0.3861386138613861,
    0.5841584158415841,
    0.06930693069306931,
    0.2772277227722773,
    0.7425742574257426,
    0.1386138613861386,
    0.8118811881188119,
    0.7722772277227723,
    0.900990099009901,
    0.4455445544554456,
    0.495049504950495,
    0.9603960396039604,
    0.0396039603960396,
    0.0396039603960396,
    0.3465346534653465,
    0.0396039603960396,
    0.7326732673267327,
    0.7425742574257426,
    0.6237623762376238,
    0.3069306930693069,
    0.6039603960396039,
    0.9702970297029703,
    0.3762376237623762,
    0.3663366336633663,
    0.4752475247524752,
    0.8712871287128713,
    0.9207920792079208,
    0.5148514851485149,
    0.0594059405940594,
    0.7722772277227723,
    0.8514851485148515,
    0.1485148514851485,
    0.8811881188118812,
    0.1782178217821782,
    0.8217821782178217,
    0.2475247524752475,
    0.09900990099009901,
    0.1782178217821782,
    0.1485148514851485,
    0.8118811881188119,
    0.5445544554455446,
    0.8217821782178217,
    0.5643564356435643,
    0.6039603960396039,
    0.06930693069306931,
    0.1188118811881188,
    0.1188118811881188,
    0.07920792079207921,
    0.2079207920792079,
    0.2475247524752475,
    0.6336633663366337,
    0.0891089108910891,
    0.2475247524752475,
    0.2673267326732673,
    0.6336633663366337,
    0.8217821782178217,
    0.0297029702970297,
    0.6039603960396039,
    0.7326732673267327,
    0.2772277227722773,
    0.5445544554455446,
    0.6039603960396039,
    0.2574257425742574,
    0.8613861386138614,
    0.2178217821782178,
    0.8118811881188119,
    0.3564356435643564,
    0.801980198019802,
    0.7920792079207921,
    0.4851485148514851,
    0.3861386138613861,
    0.7128712871287128,
    0.2475247524752475,
    0.2772277227722773,
    0.2871287128712871,
    0.009900990099009901,
    0.5544554455445545,
    0.8811881188118812,
    0.297029702970297,
    0.009900990099009901,
    0.4653465346534654,
    0.1386138613861386,
    0.8712871287128713,
    0.1782178217821782,
    0.2574257425742574,
    0.7722772277227723,
    0.405940594059406,
    0.9603960396039604,
    0.7128712871287128,
    0.8712871287128713,
    0.8514851485148515,
    0.6633663366336634,
    0.7623762376237624,
    0.3663366336633663,
    0.6336633663366337,
    0.0297029702970297,
    0.7425742574257426,
    0.3564356435643564,
    0.06930693069306931,
    0.405940594059406,
    0.1485148514851485,
    0.7821782178217822,
    0.6138613861386139,
    0.900990099009901,
    0.6831683168316832,
    0.0198019801980198,
    0.5643564356435643,
    0.2673267326732673,
    0.9801980198019802,
    0.5445544554455446,
    0.9306930693069307,
    0.1089108910891089,
    0.9306930693069307,
    0.297029702970297,
    0.5247524752475248,
    0.7425742574257426,
    0.1881188118811881,
    0.4851485148514851,
    0.9801980198019802,
    0.5544554455445545,
    0.3267326732673267,
    0.7227722772277227,
    0.5445544554455446,
    0.5247524752475248,
    0.9108910891089109,
    0.5346534653465347,
    0.297029702970297,
    0.594059405940594,
    0.6435643564356436,
    0.5841584158415841,
    0.2772277227722773,
    0.9504950495049505,
    0.2673267326732673,
    0.6633663366336634,
    0.0396039603960396,
    0.2277227722772277,
    0.1188118811881188,
    0.7326732673267327,
    0.4158415841584158,
    0.07920792079207921,
    0.1584158415841584,
    0.5247524752475248,
    0.06930693069306931,
    0.3564356435643564,
    0.7326732673267327,
    0.4356435643564356,
    0.7326732673267327,
    0.5445544554455446,
    0.07920792079207921
)

$n = $pValues.Length
$data = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $data[$i,0] = $pValues[$i]
}

$target = $ws.Range("F2:F251")
$target.Value = $data
